$d = $word.ActiveDocument

# Build the OOXML fragment for the new "Friday 15 March" section that
# gets appended after the existing final paragraph (which holds the
# _GoBack bookmark and must stay last).
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$pEmpty = "<w:p $ns><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr></w:p>"

$pHeading = "<w:p $ns><w:pPr><w:pStyle w:val=`"Heading2`"/><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" +
            "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>Friday 15 March</w:t></w:r></w:p>"

$pTime = "<w:p $ns><w:pPr><w:rPr><w:i/><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" +
         "<w:r><w:rPr><w:i/><w:lang w:val=`"en-US`"/></w:rPr><w:t>1hr</w:t></w:r></w:p>"

$bodyRuns =
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>Consolidated everyone’s progress including the deployment of the website to a domain, ‘</w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>carboncommuter.xyz</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`">’ and the technical details of running it in debug versus in production. We discussed logging and how we can implement it, and fixing some issues with regards to the </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>gitignore</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> file not updating properly. Finally, we discussed how we will integrate each other’s changes over the weekend.</w:t></w:r>"

$pBody = "<w:p $ns><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>$bodyRuns</w:p>"

$fragment = $pEmpty + $pHeading + $pTime + $pBody

# Insert the whole fragment in one shot right after the document's
# current last paragraph (which contains the _GoBack bookmark), so the
# bookmark stays attached to the new final paragraph.
$insertAt = $d.Content.End
$rng = $d.Range($insertAt, $insertAt)
[void]$rng.InsertXML($fragment)
